$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5090
$ws.Range("J3").Value = 8081
$ws.Range("L3").Value = 5482
$ws.Range("L4").Value = 1341
$ws.Range("L6").Value = 4585
$ws.Range("J7").Value = 29354
$ws.Range("L7").Value = 16822

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 323
$ws.Range("L6").Value = 286
$ws.Range("L7").Value = 1110

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 152
$ws.Range("L7").Value = 373

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 271
$ws.Range("L6").Value = 228
$ws.Range("L7").Value = 777

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 74
$ws.Range("L7").Value = 236

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 194
$ws.Range("L3").Value = 217
$ws.Range("L6").Value = 178
$ws.Range("L7").Value = 643

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 122
$ws.Range("L7").Value = 295

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L6").Value = 128
$ws.Range("L7").Value = 546
$ws.Range("L8").Value = 1110
$ws.Range("L9").Value = 98
$ws.Range("L10").Value = 108
$ws.Range("L11").Value = 273
$ws.Range("L12").Value = 37
$ws.Range("L14").Value = 88
$ws.Range("L19").Value = 459
$ws.Range("L20").Value = 415
$ws.Range("L23").Value = 183
$ws.Range("L29").Value = 926
$ws.Range("L33").Value = 777
$ws.Range("L34").Value = 98
$ws.Range("L37").Value = 643
$ws.Range("L42").Value = 549
$ws.Range("L43").Value = 122
$ws.Range("L44").Value = 116
$ws.Range("L47").Value = 114
$ws.Range("L49").Value = 85
$ws.Range("L50").Value = 87
$ws.Range("L51").Value = 209
$ws.Range("L54").Value = 360
$ws.Range("L55").Value = 167
$ws.Range("J63").Value = 228
$ws.Range("L63").Value = 50
$ws.Range("L78").Value = 214
$ws.Range("L79").Value = 443
$ws.Range("L83").Value = 373
$ws.Range("L84").Value = 164
$ws.Range("L89").Value = 245
$ws.Range("L95").Value = 236
$ws.Range("L98").Value = 93
$ws.Range("L99").Value = 295
$ws.Range("J101").Value = 29354
$ws.Range("L101").Value = 16822

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 57
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 85

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 65
$ws.Range("L3").Value = 87
$ws.Range("L6").Value = 176
$ws.Range("L7").Value = 360

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 281
$ws.Range("L3").Value = 349
$ws.Range("L6").Value = 234
$ws.Range("L7").Value = 926

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L6").Value = 128
$ws.Range("L7").Value = 459

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L3").Value = 34
$ws.Range("L7").Value = 116

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L4").Value = 8
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L3").Value = 38
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 128

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 187
$ws.Range("L7").Value = 549

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 108

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 58
$ws.Range("L7").Value = 214

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 54
$ws.Range("L7").Value = 167

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L4").Value = 17
$ws.Range("L7").Value = 183

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L6").Value = 98
$ws.Range("L7").Value = 443

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 135
$ws.Range("L7").Value = 415

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 185
$ws.Range("L3").Value = 185
$ws.Range("L4").Value = 36
$ws.Range("L7").Value = 546

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L2").Value = 31
$ws.Range("L7").Value = 98

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L2").Value = 44
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L2").Value = 30
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 83
$ws.Range("L7").Value = 273

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L2").Value = 29
$ws.Range("L7").Value = 98

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L6").Value = 67
$ws.Range("L7").Value = 245

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 66
$ws.Range("L7").Value = 209

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 122

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 37
